$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("level")

# Update the F-column JSON payloads: inject a "speed" field (and variants)
# into the mob-stats blobs that previously all shared the bare
# {"sight":1.5} / {"sight": 2.5} strings.
$ws.Range("F2").Value   = '{"sight":1.5,"speed":22.0}'
$ws.Range("F3").Value   = '{"sight":1.5,"speed":22.1}'
$ws.Range("F4").Value   = '{"sight":1.5,"speed":22.2}'
$ws.Range("F5").Value   = '{"sight":1.5,"speed":22.3}'
$ws.Range("F6").Value   = '{"sight":1.5,"speed":22.4}'
$ws.Range("F7").Value   = '{"sight":1.5,"speed":22.5}'
$ws.Range("F8").Value   = '{"sight":1.5,"speed":22.6}'
$ws.Range("F9").Value   = '{"sight":1.5,"speed":22.7}'
$ws.Range("F10").Value  = '{"sight":1.5,"speed":22.8}'
$ws.Range("F11").Value  = '{"sight":1.5,"speed":22.9}'
$ws.Range("F12").Value  = '{"sight":1.5,"speed":23.0}'
$ws.Range("F13").Value  = '{"sight":1.5,"speed":23.1}'
$ws.Range("F14").Value  = '{"sight":1.5,"speed":23.2}'
$ws.Range("F15").Value  = '{"sight":1.5,"speed":23.3}'
$ws.Range("F16").Value  = '{"sight":1.5,"speed":23.4}'
$ws.Range("F17").Value  = '{"sight":1.5,"speed":23.5}'
$ws.Range("F18").Value  = '{"sight":1.5,"speed":23.6}'
$ws.Range("F19").Value  = '{"sight":1.5,"speed":23.7}'
$ws.Range("F20").Value  = '{"sight":1.5,"speed":23.8}'
$ws.Range("F21").Value  = '{"sight":1.5,"speed":23.9}'
$ws.Range("F22").Value  = '{"sight":1.5,"speed":24}'
$ws.Range("F103").Value = '{"sight": 2.5, "speed":25}'
$ws.Range("F104").Value = '{"sight": 2.5}'

# Widen column F so the longer JSON strings are readable.
$ws.Columns.Item(6).ColumnWidth = 20.5

# Move the active selection on the frozen pane from F2 down to F3.
[void]$ws.Activate()
[void]$ws.Range("F3").Select()
